# Append nine new English/Chinese vocabulary rows to Sheet1 (rows 150-158),
# which previously held only blank placeholder cells.
#
# The writes are ordered to reproduce the exact shared-string table layout
# of the target workbook: most rows were filled in column-A-then-column-B
# order, but the translation (column B) for "in that regard" and
# "falling edge" was entered before their English term (column A) - and
# those two English terms were typed in last, after every other row was
# already complete.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A150").Value = "configuration"
$ws.Range("B150").Value = "配置"
$ws.Range("B151").Value = "在这方面"
$ws.Range("A152").Value = "appropriately"
$ws.Range("B152").Value = "适当"
$ws.Range("A153").Value = "configuraion"
$ws.Range("B153").Value = "结构"
$ws.Range("A154").Value = "individually"
$ws.Range("B154").Value = "个别地"
$ws.Range("A155").Value = "interconnect"
$ws.Range("B155").Value = "互连"
$ws.Range("A156").Value = "throughtput"
$ws.Range("B156").Value = "输出通量"
$ws.Range("B157").Value = "下降沿"
$ws.Range("A158").Value = "ramble"
$ws.Range("B158").Value = "漫无目的"
$ws.Range("A151").Value = "in that regard"
$ws.Range("A157").Value = "falling edge"

$ws.Range("B160").Select()
